$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 108: new record for "2018 Second Primary"
$ws.Range("A107:N107").Copy()
$ws.Range("A108:N108").PasteSpecial(-4122)
$ws.Range("D108").Clear()
$ws.Range("I108").Clear()
$ws.Range("J108").Clear()
$ws.Range("N108").Clear()

$ws.Range("A108").Value = 128
$ws.Range("B108").Value = 43234.9993055556
$ws.Range("C108").Value = 43227.9993055556
$ws.Range("E108").Value = "2018 Second Primary"
$ws.Range("F108").Value = 1
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 0
$ws.Range("K108").Value = 43193
$ws.Range("L108").Value = 1
$ws.Range("M108").Value = 1

# Row 109: new record for "2018 Second Biannual"
$ws.Range("A107:N107").Copy()
$ws.Range("A109:N109").PasteSpecial(-4122)
$ws.Range("D109").Clear()
$ws.Range("I109").Clear()
$ws.Range("J109").Clear()
$ws.Range("N109").Clear()

$ws.Range("A109").Value = 129
$ws.Range("B109").Value = 43381.9993055556
$ws.Range("C109").Value = 43374.9993055556
$ws.Range("E109").Value = "2018 Second Biannual"
$ws.Range("F109").Value = 1
$ws.Range("G109").Value = 3
$ws.Range("H109").Value = 0
$ws.Range("K109").Value = 43193
$ws.Range("L109").Value = 1
$ws.Range("M109").Value = 1
